# HMP2360_r0_t0.xlsx - kinetics1 sheet: split the single "order" column into
# two columns: "substrate order" and "product order" (fix per commit message
# "Fixed issue with order column, should be 2 columns: substrate order and
# product order").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kinetics1")

# Insert a brand-new column before D - this shifts the old D..K
# (promiscuous..comments) one place to the right, to E..L, and leaves a
# blank column D for the new "product order" data.
$ws.Columns("D").Insert()

# C1 held "order" - rename it to "substrate order". D1 is the new column,
# give it the "product order" header.
$ws.Range("C1").Value = "substrate order"
$ws.Range("D1").Value = "product order"

# The previous single "order" value applied to both the substrate and the
# product side, so seed the new "product order" column with a copy of
# whatever is in "substrate order" for each data row (rows 2-14).
for ($r = 2; $r -le 14; $r++) {
    $substrateOrder = $ws.Cells.Item($r, 3).Value2
    if ($substrateOrder -ne $null) {
        $ws.Cells.Item($r, 4).Value = $substrateOrder
    }
}

# Widen the two text-heavy columns now that they hold real header/content
# text, and tighten the row height slightly to match the edited layout.
$ws.Columns("C").ColumnWidth = 22.5
$ws.Columns("D").ColumnWidth = 23.5
$ws.Rows("1:14").RowHeight = 13.8

# Make kinetics1 the active sheet/tab, with D2:D7 (the newly-populated
# "product order" values) selected.
$ws.Activate()
$ws.Range("D2:D7").Select()
